$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list1_sheet1"

# --- Add the two new sheets after the existing one, in order ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "list2_sheet2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "list3_sheet3"

# --- Column widths (same layout as sheet1) ---
foreach ($ws in @($ws2, $ws3)) {
    $ws.Columns.Item(1).ColumnWidth = 13.36328125
    $ws.Columns.Item(2).ColumnWidth = 12.7265625
    $ws.Columns.Item(3).ColumnWidth = 21.6328125
    $ws.Columns.Item(4).ColumnWidth = 12.453125
    $ws.Columns.Item(5).ColumnWidth = 10.08984375
}

# --- Header rows (bold style, matching sheet1's header) ---
$headers = @("First Name", "Last Name", "Email", "City", "Country")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 1
    $ws2.Cells.Item(1, $col).Value = $headers[$i]
    $ws3.Cells.Item(1, $col).Value = $headers[$i]
}
$ws2.Range("A1:E1").Font.Bold = $true
$ws3.Range("A1:E1").Font.Bold = $true

# --- list2_sheet2 data rows ---
$ws2.Cells.Item(2, 1).Value = "Macy"
$ws2.Cells.Item(2, 2).Value = "Barker"
$ws2.Cells.Item(2, 3).Value = "macy@gmail.com"
$ws2.Cells.Item(2, 4).Value = "New Jersey"
$ws2.Cells.Item(2, 5).Value = "USA"

$ws2.Cells.Item(3, 1).Value = "Casper"
$ws2.Cells.Item(3, 2).Value = "Pitts"
$ws2.Cells.Item(3, 3).Value = "casper@gmail.com"
$ws2.Cells.Item(3, 4).Value = "Seattle"
$ws2.Cells.Item(3, 5).Value = "USA"

# --- list3_sheet3 data rows ---
$ws3.Cells.Item(2, 1).Value = "Aleena"
$ws3.Cells.Item(2, 2).Value = "Cobb"
$ws3.Cells.Item(2, 3).Value = "aleena@gmail.com"
$ws3.Cells.Item(2, 4).Value = "Chicago"
$ws3.Cells.Item(2, 5).Value = "USA"

$ws3.Cells.Item(3, 1).Value = "Murray"
$ws3.Cells.Item(3, 2).Value = "Smart"
$ws3.Cells.Item(3, 3).Value = "murray@gmail.com"
$ws3.Cells.Item(3, 4).Value = "Phoenix"
$ws3.Cells.Item(3, 5).Value = "USA"

# --- Hyperlinks on the Email column ---
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:macy@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:casper@gmail.com")

$ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:aleena@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:murray@gmail.com")

# Blank rows below the data keep the "Email column" hyperlink styling,
# matching the two trailing formatted-but-empty rows under each table.
$ws2.Range("C4").Style = "Hyperlink"
$ws2.Range("C5").Style = "Hyperlink"
$ws3.Range("C4").Style = "Hyperlink"
$ws3.Range("C5").Style = "Hyperlink"

# --- Page setup (portrait, matching sheet1) ---
$ws2.PageSetup.Orientation = 1
$ws3.PageSetup.Orientation = 1

# --- Selections matching the target state ---
[void]$ws2.Range("D3").Select()
[void]$ws3.Range("C4").Select()

# list3_sheet3 is the active tab when the file is saved
[void]$ws3.Select()

Write-Output "Done"
